# Scheduled-runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# figures across several sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 136
$ws.Range("I28").Value = 14.666667
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 14.666667
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = 470.333333
$ws.Range("N28").Value = -1470

$ws.Range("H62").Value = 2617.5
$ws.Range("I62").Value = 2561
$ws.Range("K62").Value = 2561
$ws.Range("M62").Value = -1937

$ws.Range("H65").Value = 2617.5
$ws.Range("I65").Value = 2561
$ws.Range("K65").Value = 12805
$ws.Range("M65").Value = -9685

$ws.Range("H112").Value = 1793088.1
$ws.Range("J112").Value = 1822474.9
$ws.Range("L112").Value = 5467424.699999999
$ws.Range("N112").Value = -5469640.699999999

$ws.Range("H129").Value = 854.3200000000001
$ws.Range("I129").Value = 800
$ws.Range("J129").Value = 855.4286
$ws.Range("K129").Value = 2400
$ws.Range("L129").Value = 2566.2858
$ws.Range("M129").Value = 2600
$ws.Range("N129").Value = -12566.2858

$ws.Range("H138").Value = 1904.5111
$ws.Range("I138").Value = 563.73914
$ws.Range("J138").Value = 3306.2273
$ws.Range("K138").Value = 1691.21742
$ws.Range("L138").Value = 9918.6819
$ws.Range("M138").Value = 3448.78258
$ws.Range("N138").Value = -20198.6819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20381.982
$ws.Range("I32").Value = 21079.04
$ws.Range("J32").Value = 8299.666999999999
$ws.Range("K32").Value = 21079.04
$ws.Range("L32").Value = 8299.666999999999
$ws.Range("M32").Value = -20792.04
$ws.Range("N32").Value = -8873.666999999999

$ws.Range("H61").Value = 1564.1132
$ws.Range("J61").Value = 3888.2856
$ws.Range("L61").Value = 3888.2856
$ws.Range("N61").Value = -4312.2856

$ws.Range("H136").Value = 1564.1132
$ws.Range("J136").Value = 3888.2856
$ws.Range("L136").Value = 11664.8568
$ws.Range("N136").Value = -16764.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2014.9333
$ws.Range("I86").Value = 1768.6666
$ws.Range("K86").Value = 1768.6666
$ws.Range("M86").Value = -645.6666

$ws.Range("H89").Value = 2014.9333
$ws.Range("I89").Value = 1768.6666
$ws.Range("K89").Value = 8843.333000000001
$ws.Range("M89").Value = -3227.333000000001

$ws.Range("H99").Value = 1680.75
$ws.Range("J99").Value = 1400
$ws.Range("L99").Value = 1400
$ws.Range("N99").Value = -4396

$ws.Range("H134").Value = 26729.635
$ws.Range("I134").Value = 28742.5
$ws.Range("J134").Value = 1233.3334
$ws.Range("K134").Value = 86227.5
$ws.Range("L134").Value = 3700.0002
$ws.Range("M134").Value = -83692.5
$ws.Range("N134").Value = -8770.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 255
$ws.Range("I22").Value = 255
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 255
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = 95

$ws.Range("H58").Value = 19783.852
$ws.Range("I58").Value = 1446.2667
$ws.Range("J58").Value = 42705.832
$ws.Range("K58").Value = 1446.2667
$ws.Range("L58").Value = 42705.832
$ws.Range("M58").Value = -1243.2667
$ws.Range("N58").Value = -43111.832

$ws.Range("H62").Value = 6602.4
$ws.Range("J62").Value = 7003
$ws.Range("L62").Value = 7003
$ws.Range("N62").Value = -8251

$ws.Range("H65").Value = 6602.4
$ws.Range("J65").Value = 7003
$ws.Range("L65").Value = 35015
$ws.Range("N65").Value = -41255

$ws.Range("H99").Value = 14230807
$ws.Range("I99").Value = 2978605
$ws.Range("J99").Value = 38466320
$ws.Range("K99").Value = 2978605
$ws.Range("L99").Value = 38466320
$ws.Range("M99").Value = -2977107
$ws.Range("N99").Value = -38469316

$ws.Range("H126").Value = 14230807
$ws.Range("I126").Value = 2978605
$ws.Range("J126").Value = 38466320
$ws.Range("K126").Value = 8935815
$ws.Range("L126").Value = 115398960
$ws.Range("M126").Value = -8933345
$ws.Range("N126").Value = -115403900

$ws.Range("H136").Value = 19783.852
$ws.Range("I136").Value = 1446.2667
$ws.Range("J136").Value = 42705.832
$ws.Range("K136").Value = 4338.800099999999
$ws.Range("L136").Value = 128117.496
$ws.Range("M136").Value = -1788.800099999999
$ws.Range("N136").Value = -133217.496

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 241.71428
$ws.Range("I14").Value = 241.71428
$ws.Range("K14").Value = 725.14284
$ws.Range("M14").Value = -552.14284

$ws.Range("H131").Value = 760.37
$ws.Range("I131").Value = 360
$ws.Range("J131").Value = 777.05206
$ws.Range("K131").Value = 1080
$ws.Range("L131").Value = 2331.15618
$ws.Range("M131").Value = 3960
$ws.Range("N131").Value = -12411.15618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5686300
$ws.Range("I70").Value = 4750
$ws.Range("J70").Value = 6948866.5
$ws.Range("K70").Value = 4750
$ws.Range("L70").Value = 6948866.5
$ws.Range("M70").Value = -4480
$ws.Range("N70").Value = -6949406.5

$ws.Range("H73").Value = 5686300
$ws.Range("I73").Value = 4750
$ws.Range("J73").Value = 6948866.5
$ws.Range("K73").Value = 4750
$ws.Range("L73").Value = 6948866.5
$ws.Range("M73").Value = -3814
$ws.Range("N73").Value = -6950738.5

$ws.Range("H80").Value = 3712.2222
$ws.Range("I80").Value = 3480.125
$ws.Range("K80").Value = 3480.125
$ws.Range("M80").Value = -2482.125

$ws.Range("H83").Value = 3712.2222
$ws.Range("I83").Value = 3480.125
$ws.Range("K83").Value = 17400.625
$ws.Range("M83").Value = -12408.625

$ws.Range("H126").Value = 4790.909
$ws.Range("I126").Value = 3460
$ws.Range("J126").Value = 7642.857
$ws.Range("K126").Value = 10380
$ws.Range("L126").Value = 22928.571
$ws.Range("M126").Value = -7910
$ws.Range("N126").Value = -27868.571

$ws.Range("H132").Value = 103519.266
$ws.Range("I132").Value = 94253.91
$ws.Range("J132").Value = 128999
$ws.Range("K132").Value = 282761.73
$ws.Range("L132").Value = 386997
$ws.Range("M132").Value = -280231.73
$ws.Range("N132").Value = -392057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 205.5625
$ws.Range("I55").Value = 155.6
$ws.Range("J55").Value = 228.27272
$ws.Range("K55").Value = 155.6
$ws.Range("L55").Value = 228.27272
$ws.Range("M55").Value = 17.40000000000001
$ws.Range("N55").Value = -574.2727199999999

$ws.Range("H136").Value = 21744.2
$ws.Range("I136").Value = 36900.355
$ws.Range("J136").Value = 2454.5454
$ws.Range("K136").Value = 110701.065
$ws.Range("L136").Value = 7363.6362
$ws.Range("M136").Value = -108151.065
$ws.Range("N136").Value = -12463.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 801.41174
$ws.Range("I100").Value = 509.53845
$ws.Range("K100").Value = 1019.0769
$ws.Range("M100").Value = -478.0769

Write-Host "Applied scheduled-runner updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR."
